$d = $word.ActiveDocument
$n = $d.Paragraphs.Count

# The site-footer block to drop consists of three consecutive paragraphs:
#   1. a blank paragraph
#   2. "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3. "© 2020 . Contact: ... Creative Commons Attribution"
# Locate paragraph 2 by its text, then remove it together with its immediate
# neighbours (the blank paragraph before it and the copyright paragraph after
# it), leaving the "LOT2007: ..." paragraph and everything that follows the
# footer block untouched.
$verIndex = 0
for ($i = 1; $i -le $n; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*Ver no Jupiter Salvar em pdf Salvar em docx*") {
        $verIndex = $i
        break
    }
}

if ($verIndex -gt 0) {
    $startPara = $d.Paragraphs.Item($verIndex - 1)
    $endPara = $d.Paragraphs.Item($verIndex + 1)
    $r = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $r.Delete()
}
